$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Расширенные атрибуты увидеть не удалось - гостевому пользователю отказано в доступе.*") {
        $p.Range.Delete()
        break
    }
}
